# Update data models and national cases
# Fill in the previously-empty rows 202-204 on "Hoja1" with the next three
# days of data (Dia, Casos Activos, Casos Confirmados), and move the active
# selection to C200 (matches where the sheet view now points after the
# new rows were entered).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 202: day 29
$ws.Range("A202").Value = 29
$ws.Range("B202").Value = 1034
$ws.Range("C202").Value = 18403

# Row 203: day 30
$ws.Range("A203").Value = 30
$ws.Range("B203").Value = 1014
$ws.Range("C203").Value = 18511

# Row 204: day 1 (new month)
$ws.Range("A204").Value = 1
$ws.Range("B204").Value = 1033
$ws.Range("C204").Value = 18628

# Move / record the active selection as it ends up after the data entry
$ws.Range("C200").Select()
